# Generate Report for Handoff
#
# The localized docs just got queued for handoff again, so refresh the
# status text + the "xliff generated" / handoff timestamps on every
# language sheet, and let the Status column shrink to fit the new
# (shorter) text.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet: status appears twice (zh-cn + de-de columns) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-12 05:02:26"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-12 05:02:21"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-12 05:02:26"

# --- Shrink the Status column now that the text is shorter ---
# (ColumnWidth is expressed in characters; 16.3333... is the closest
# setting to the fitted width produced by Excel for this text/font.)
$fitWidth = 16.333333333333336
$overview.Columns.Item(5).ColumnWidth = $fitWidth
$overview.Columns.Item(6).ColumnWidth = $fitWidth
$zhcn.Columns.Item(3).ColumnWidth = $fitWidth
$dede.Columns.Item(3).ColumnWidth = $fitWidth
